# SCD0016-032 - Searching profiling nasabah & memiliki sales kelolaan
# Update Excel SCD0011 until SCD0016
#
# This script reproduces, via Excel COM automation, the edits captured in the
# canonical OOXML diff:
#   - sheet renamed SCD0258 -> SCD0016
#   - TC_ID cell (B2) updated from "DGS-273" to "SCD0016-032"
#   - whole-sheet re-alignment: horizontal=left, vertical=center applied
#     to every populated cell (this is what collapsed/renumbered cellXfs)
#   - a handful of now-redundant blank formatted cells were cleared
#   - active selection moved to B3, view scrolled back to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the worksheet -------------------------------------------------
$ws.Name = "SCD0016"

# --- update the TC_ID value in B2 ------------------------------------------
$ws.Range("B2").Value = "SCD0016-032"
# B2 picks up the same font used by the rest of column B's data rows (Arial 10)
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10

# --- re-align every populated cell (horizontal=left, vertical=center) ------
$alignRanges = @(
    "A1:N1",
    "A2:P2", "R2:S2",
    "A3:I3", "O3:P3",
    "A4:I4", "N4:Q4",
    "B5", "F5:G5",
    "B6", "F6:G6"
)
foreach ($addr in $alignRanges) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4131   # xlLeft
    $rng.VerticalAlignment = -4108     # xlCenter
}

# --- clear cells that become redundant blanks after the re-alignment -------
$clearCells = @("I2", "J2", "K2", "P2", "R2", "S2", "A3", "H3", "A4", "H4")
foreach ($addr in $clearCells) {
    $ws.Range($addr).Clear()
}

# --- restore the view: scroll to A2, select B3 ------------------------------
$ws.Range("A2").Select()
$ws.Range("B3").Select()
